$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.999.03'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '3.428.54'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '412.48'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.81'
$ws.Range("E6").Value = '  +0.65%  '
$ws.Range("E7").Value = '  +1.43%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.729'
$ws.Range("E9").Value = '  -2.62%  '
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.60'
$ws.Range("E11").Value = '  +1.34%  '
$ws.Range("E12").Value = '  +4.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000217'
$ws.Range("E13").Value = '  +7.27%  '
$ws.Range("D14").Value = '3.969.69'
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.13'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = '3.433.24'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("E18").Value = '  +4.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.57'
$ws.Range("E19").Value = '  -1.26%  '
$ws.Range("D20").Value = '61.945.16'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '491.91'
$ws.Range("E21").Value = '  +21.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '92.62'
$ws.Range("E22").Value = '  +2.51%  '
$ws.Range("E23").Value = '  +4.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.64'
$ws.Range("E24").Value = '  +1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.39'
$ws.Range("E25").Value = '  +4.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '34.63'
$ws.Range("E26").Value = '  +5.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.15'
$ws.Range("E27").Value = '  +7.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.80'
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.75'
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.09'
$ws.Range("E31").Value = '  +2.15%  '
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("E33").Value = '  -3.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.11'
$ws.Range("E34").Value = '  -4.37%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.17'
$ws.Range("E36").Value = '  +9.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0497'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.47'
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '151.21'
$ws.Range("E40").Value = '  +7.21%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.329'
$ws.Range("E41").Value = '  +4.34%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  +3.41%  '
$ws.Range("E43").Value = '  +9.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.98'
$ws.Range("E44").Value = '  +2.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.71'
$ws.Range("E45").Value = '  +14.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.28'
$ws.Range("E46").Value = '  +6.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.36'
$ws.Range("E47").Value = '  +21.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.69'
$ws.Range("E48").Value = '  -0.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.89'
$ws.Range("E49").Value = '  +4.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '118.08'
$ws.Range("E50").Value = '  +23.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.145'
$ws.Range("E51").Value = '  +13.60%  '
